# Update the Drools decision-table worksheet ("Sheet1") with the refactored
# rule-table content: the ProductDetails/UserDetails bindings are now bound
# to 'product'/'user' variables, the conditions use explicit boolean
# comparisons, the action calls setDisplayProduct() on the bound variable,
# and the ruleset import cell / trailing "Variables" row are cleaned up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ruleset name cell
$ws.Range("C1").Value = "poc.decisiontable.fuse.brms.displayproduct"

# RuleTable parameter bindings (row 7)
$ws.Range("B7").Value = "product:ProductDetails"
$ws.Range("C7").Value = "user:UserDetails"
$ws.Range("D7").ClearContents()

# RuleTable condition/action templates (row 8)
$ws.Range("B8").Value = 'productType=="$param"'
$ws.Range("C8").Value = 'userRole=="$param"'
$ws.Range("D8").Value = 'product.setDisplayProduct("$param");'

# Remove the now-unused "Variables" helper row at the bottom
$ws.Range("B25").ClearContents()
$ws.Range("C25").ClearContents()

# Restore the selected cell as left by the author
$ws.Range("C9").Select()
